$wb = $excel.ActiveWorkbook

# ---- Phase 1: rename all sheets to temporary unique names ----
$ws = $wb.Worksheets.Item(1)
$ws.Name = "__tmp_sheet_1__"
$ws = $wb.Worksheets.Item(2)
$ws.Name = "__tmp_sheet_2__"
$ws = $wb.Worksheets.Item(3)
$ws.Name = "__tmp_sheet_3__"
$ws = $wb.Worksheets.Item(4)
$ws.Name = "__tmp_sheet_4__"
$ws = $wb.Worksheets.Item(5)
$ws.Name = "__tmp_sheet_5__"
$ws = $wb.Worksheets.Item(6)
$ws.Name = "__tmp_sheet_6__"
$ws = $wb.Worksheets.Item(7)
$ws.Name = "__tmp_sheet_7__"
$ws = $wb.Worksheets.Item(8)
$ws.Name = "__tmp_sheet_8__"
$ws = $wb.Worksheets.Item(9)
$ws.Name = "__tmp_sheet_9__"

# ---- Sheet 1: summ19 -> summ5 ----
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ5"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.4139183456802575"
$ws.Range("C2").Value = [double]"0.5971320392702484"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2693890147617579"
$ws.Range("C3").Value = [double]"0.04959447246122681"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.5974541915940856"
$ws.Range("C4").Value = [double]"6.262914010172245e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3771080200757158"
$ws.Range("C5").Value = [double]"0.0008832734668545563"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3953898511536404"
$ws.Range("C6").Value = [double]"0.001836441590167969"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.01059232938917822"
$ws.Range("C7").Value = [double]"0.8411263821233601"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006171806460343455"
$ws.Range("C8").Value = [double]"2.111378898907939e-94"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01405391521876646"
$ws.Range("C9").Value = [double]"2.282773732821315e-07"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.2675592382573362"
$ws.Range("C10").Value = [double]"0.0001515461199274023"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.3177176312531071"
$ws.Range("C11").Value = [double]"0.004457876594439292"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.3295635194820375"
$ws.Range("C12").Value = [double]"0.02097263070519777"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.686829979362718e-05"
$ws.Range("C13").Value = [double]"3.258716314149012e-05"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.349679798559234e-08"
$ws.Range("C14").Value = [double]"0.2171136939498692"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.006771838159081624"
$ws.Range("C15").Value = [double]"0.7098478012521685"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04303548615439926"
$ws.Range("C16").Value = [double]"0.004259804537852114"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.371349882427615"
$ws.Range("C17").Value = [double]"5.597648384830495e-11"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.004846261554518002"
$ws.Range("C18").Value = [double]"0.1526596498317285"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.0022363436512272"
$ws.Range("C19").Value = [double]"0.6294524800398158"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.047317960025906"
$ws.Range("C20").Value = [double]"0.005268808906685529"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"-0.05896959121999019"
$ws.Range("C21").Value = [double]"0.9003696156395623"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.089383791511164"
$ws.Range("C22").Value = [double]"0.05243155735769115"

# ---- Sheet 2: summ3 -> summ2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ2"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.4645419887824939"
$ws.Range("C2").Value = [double]"0.5519852278546833"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2985259780401989"
$ws.Range("C3").Value = [double]"0.0279302751752839"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.6080325713543353"
$ws.Range("C4").Value = [double]"3.096004515925621e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3693198325362771"
$ws.Range("C5").Value = [double]"0.001052343605420668"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.4116222924900952"
$ws.Range("C6").Value = [double]"0.001153598319367286"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.009536805181946615"
$ws.Range("C7").Value = [double]"0.8548686084042613"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006197147926678925"
$ws.Range("C8").Value = [double]"2.087742937908902e-94"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01530840038012003"
$ws.Range("C9").Value = [double]"1.634620794408704e-08"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.2499744251970854"
$ws.Range("C10").Value = [double]"0.0004069831809133131"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2277530229151328"
$ws.Range("C11").Value = [double]"0.03937943029535874"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.2253409125604217"
$ws.Range("C12").Value = [double]"0.1118440840430332"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-3.357742990668401e-05"
$ws.Range("C13").Value = [double]"2.332401151013656e-07"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-2.759244602839995e-09"
$ws.Range("C14").Value = [double]"0.8021173896664231"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"0.001452593877287554"
$ws.Range("C15").Value = [double]"0.9364873547999122"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04924128952734591"
$ws.Range("C16").Value = [double]"0.001080284470338091"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.028795332441059"
$ws.Range("C17").Value = [double]"1.807511841723751e-08"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.004842544450111341"
$ws.Range("C18").Value = [double]"0.1545372536652748"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.003619515614914427"
$ws.Range("C19").Value = [double]"0.4360097784182849"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.188005366786903"
$ws.Range("C20").Value = [double]"0.001578396828830057"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"-0.06109426125563058"
$ws.Range("C21").Value = [double]"0.8973784539830738"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.052310254965676"
$ws.Range("C22").Value = [double]"0.05806942530978309"

# ---- Sheet 3: summ9 -> summ4 ----
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ4"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.5290593970035119"
$ws.Range("C2").Value = [double]"0.5046298266903471"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2004904132664882"
$ws.Range("C3").Value = [double]"0.1395233073930607"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.5999191433465174"
$ws.Range("C4").Value = [double]"4.879625763572893e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.4103346033231992"
$ws.Range("C5").Value = [double]"0.0002765618916560832"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.4325801621696042"
$ws.Range("C6").Value = [double]"0.0006490544505453466"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.02777158769730377"
$ws.Range("C7").Value = [double]"0.594929531189514"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006134664470336997"
$ws.Range("C8").Value = [double]"1.396723180968373e-93"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01425811348731194"
$ws.Range("C9").Value = [double]"1.367556206129435e-07"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.3203695845048209"
$ws.Range("C10").Value = [double]"5.902882328046963e-06"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2567721654052486"
$ws.Range("C11").Value = [double]"0.01994490204170757"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.2642289829050251"
$ws.Range("C12").Value = [double]"0.06253271836449804"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.385381581486945e-05"
$ws.Range("C13").Value = [double]"0.0002433675292828677"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.558329489543996e-08"
$ws.Range("C14").Value = [double]"0.1613642414695359"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.001756439001064855"
$ws.Range("C15").Value = [double]"0.9237814353030175"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.05593951044450174"
$ws.Range("C16").Value = [double]"0.0002508322615645413"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-1.91195671030818"
$ws.Range("C17").Value = [double]"1.117772365040974e-07"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.006801385899461691"
$ws.Range("C18").Value = [double]"0.04699357506278973"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.00309295794445747"
$ws.Range("C19").Value = [double]"0.5104347079365219"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.000384291431053"
$ws.Range("C20").Value = [double]"0.007794704472394028"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"0.04470789906740787"
$ws.Range("C21").Value = [double]"0.9246790585531731"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-0.9117057623035423"
$ws.Range("C22").Value = [double]"0.1066589828600744"

# ---- Sheet 4: summ2 -> summ0 ----
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ0"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.5753996282878451"
$ws.Range("C2").Value = [double]"0.4610227779851999"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2755764784587694"
$ws.Range("C3").Value = [double]"0.04347813331945452"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.5723592521343855"
$ws.Range("C4").Value = [double]"2.013310944857316e-07"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3453316394622928"
$ws.Range("C5").Value = [double]"0.002280514025128979"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3941852177628999"
$ws.Range("C6").Value = [double]"0.001967350957203715"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.01869932391269433"
$ws.Range("C7").Value = [double]"0.7208096423594972"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006259235685857998"
$ws.Range("C8").Value = [double]"3.571394118714741e-96"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.014885872401658"
$ws.Range("C9").Value = [double]"4.663466311695981e-08"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.2692416338462839"
$ws.Range("C10").Value = [double]"0.0001361438305075974"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.336394963903798"
$ws.Range("C11").Value = [double]"0.002915466708269253"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.3610518570279483"
$ws.Range("C12").Value = [double]"0.01206284040097568"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.822927725100847e-05"
$ws.Range("C13").Value = [double]"1.418556192722522e-05"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.280833481487032e-08"
$ws.Range("C14").Value = [double]"0.2434939637169536"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.002090275289789069"
$ws.Range("C15").Value = [double]"0.9086641934455205"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04752123859364833"
$ws.Range("C16").Value = [double]"0.001536345274111043"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-1.998359153294701"
$ws.Range("C17").Value = [double]"3.721359712829598e-08"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.006318410216525099"
$ws.Range("C18").Value = [double]"0.06259415133641066"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.003228385365745796"
$ws.Range("C19").Value = [double]"0.4894752021862241"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.305194573739011"
$ws.Range("C20").Value = [double]"0.0005032240775456324"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"0.2582498581566304"
$ws.Range("C21").Value = [double]"0.5836095757305375"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.124265212833533"
$ws.Range("C22").Value = [double]"0.0423312412622656"

# ---- Sheet 5: summ8 -> summ10 ----
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ10"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.3941535352021986"
$ws.Range("C2").Value = [double]"0.6130279583009326"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2729824927401734"
$ws.Range("C3").Value = [double]"0.04551986117735479"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.6273808177926924"
$ws.Range("C4").Value = [double]"1.217879221132238e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3619649525969827"
$ws.Range("C5").Value = [double]"0.001398857390520084"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3587022057792132"
$ws.Range("C6").Value = [double]"0.004463324909754476"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.01476086713010644"
$ws.Range("C7").Value = [double]"0.7784992272829664"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.000617003902895816"
$ws.Range("C8").Value = [double]"2.499395160545216e-94"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01540845754338325"
$ws.Range("C9").Value = [double]"1.420942051893959e-08"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.3050137828381532"
$ws.Range("C10").Value = [double]"1.721753022291158e-05"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2174412194983987"
$ws.Range("C11").Value = [double]"0.05030207753608983"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.2694449612417538"
$ws.Range("C12").Value = [double]"0.05967709774214076"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.827400305922828e-05"
$ws.Range("C13").Value = [double]"1.257109069977279e-05"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.256305077779212e-08"
$ws.Range("C14").Value = [double]"0.2546449609038404"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.005030667546445898"
$ws.Range("C15").Value = [double]"0.7813762197057763"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04528683510850862"
$ws.Range("C16").Value = [double]"0.002863181973138376"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.136418001113757"
$ws.Range("C17").Value = [double]"3.1167525068804e-09"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.007754405149491486"
$ws.Range("C18").Value = [double]"0.02310810329703295"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.003525802661747429"
$ws.Range("C19").Value = [double]"0.4490467007843965"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.167288809479434"
$ws.Range("C20").Value = [double]"0.001821329990743903"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"-0.02251224825224352"
$ws.Range("C21").Value = [double]"0.9616152872967012"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-0.8826049950535083"
$ws.Range("C22").Value = [double]"0.1124699741530894"

# ---- Sheet 6: summ4 -> summ8 ----
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ8"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.1929283246839266"
$ws.Range("C2").Value = [double]"0.806547682633473"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.3140528528185846"
$ws.Range("C3").Value = [double]"0.02022916559432216"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.6326225203578596"
$ws.Range("C4").Value = [double]"7.035172016310586e-09"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3906245636693332"
$ws.Range("C5").Value = [double]"0.0004965275642033442"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3845854413583799"
$ws.Range("C6").Value = [double]"0.002278846191356602"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.002673908268546765"
$ws.Range("C7").Value = [double]"0.9586042377715109"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006168129863991184"
$ws.Range("C8").Value = [double]"5.434525613892257e-95"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.0164412850362595"
$ws.Range("C9").Value = [double]"1.572668794504059e-09"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.2771800268663574"
$ws.Range("C10").Value = [double]"9.362220308626699e-05"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.3598375742177406"
$ws.Range("C11").Value = [double]"0.001358878795491342"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.3285364293192051"
$ws.Range("C12").Value = [double]"0.02191558151704509"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.941292785968473e-05"
$ws.Range("C13").Value = [double]"6.008910176407829e-06"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-5.989200229909213e-09"
$ws.Range("C14").Value = [double]"0.5839574527278224"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"0.0002629871781154096"
$ws.Range("C15").Value = [double]"0.9884965689385996"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.03989195928714426"
$ws.Range("C16").Value = [double]"0.007803200015280198"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.179525703753995"
$ws.Range("C17").Value = [double]"1.571396508221709e-09"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.009484974891652953"
$ws.Range("C18").Value = [double]"0.005532701058199699"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.005482757838224049"
$ws.Range("C19").Value = [double]"0.244876213453289"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.138602381924326"
$ws.Range("C20").Value = [double]"0.002160141858142881"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"-0.02327525858387971"
$ws.Range("C21").Value = [double]"0.9604084557298057"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-0.9109249829519376"
$ws.Range("C22").Value = [double]"0.1012630077969169"

# ---- Sheet 7: summ0 -> summ1 ----
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ1"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.237414392156268"
$ws.Range("C2").Value = [double]"0.7648857523929591"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2304694934566151"
$ws.Range("C3").Value = [double]"0.08918973288361121"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.6042459378641299"
$ws.Range("C4").Value = [double]"3.680759946596099e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3708089690223252"
$ws.Range("C5").Value = [double]"0.00104642665087565"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.324745812549009"
$ws.Range("C6").Value = [double]"0.01012077025950743"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.01035968094143993"
$ws.Range("C7").Value = [double]"0.841197390239778"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006501477941041793"
$ws.Range("C8").Value = [double]"1.198900785025898e-101"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01524791227834781"
$ws.Range("C9").Value = [double]"2.234085329594712e-08"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.2666048276013943"
$ws.Range("C10").Value = [double]"0.0001717444567731661"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2402140440346167"
$ws.Range("C11").Value = [double]"0.03133308049136405"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.2518566668338394"
$ws.Range("C12").Value = [double]"0.0782348616506504"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.915722240336696e-05"
$ws.Range("C13").Value = [double]"7.40124003555545e-06"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-7.683263437468269e-09"
$ws.Range("C14").Value = [double]"0.4847938376525756"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.00599240109814937"
$ws.Range("C15").Value = [double]"0.742653846210021"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04860209881459097"
$ws.Range("C16").Value = [double]"0.001304195545828684"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.051707143930902"
$ws.Range("C17").Value = [double]"1.581017391131742e-08"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.007019021117688432"
$ws.Range("C18").Value = [double]"0.0398355975401621"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.005840712041841037"
$ws.Range("C19").Value = [double]"0.2174777908196336"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.181096820299502"
$ws.Range("C20").Value = [double]"0.00173890251338312"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"0.1286232149354404"
$ws.Range("C21").Value = [double]"0.7869119027130342"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.01036582790634"
$ws.Range("C22").Value = [double]"0.07273791760873573"

# ---- Sheet 8: summ30 -> summ20 ----
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ20"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"0.02468770298177458"
$ws.Range("C2").Value = [double]"0.9753530250992769"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.21529719844816"
$ws.Range("C3").Value = [double]"0.1201978748847771"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.5968273845918851"
$ws.Range("C4").Value = [double]"7.740282438269154e-08"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3424330897902683"
$ws.Range("C5").Value = [double]"0.002701166312105048"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3720232556205244"
$ws.Range("C6").Value = [double]"0.00364437070996027"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.0288547644420757"
$ws.Range("C7").Value = [double]"0.5945369846597194"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.00064394537380495"
$ws.Range("C8").Value = [double]"1.768995044276255e-100"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01489917066174146"
$ws.Range("C9").Value = [double]"4.628985411450464e-08"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.236289707263"
$ws.Range("C10").Value = [double]"0.0008749475259994023"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2387595829366609"
$ws.Range("C11").Value = [double]"0.03217991642545891"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.2625137983789679"
$ws.Range("C12").Value = [double]"0.0673651510954188"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-3.187744560290437e-05"
$ws.Range("C13").Value = [double]"1.028196034683109e-06"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.632566644654832e-08"
$ws.Range("C14").Value = [double]"0.1385659001736696"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"0.005144240911122681"
$ws.Range("C15").Value = [double]"0.7781920245882151"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.0305818193266104"
$ws.Range("C16").Value = [double]"0.04345443675168315"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.171584006775247"
$ws.Range("C17").Value = [double]"2.404836914044022e-09"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.006031374557724413"
$ws.Range("C18").Value = [double]"0.07991984319064926"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.002272009040255304"
$ws.Range("C19").Value = [double]"0.6291165615397134"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.374624511450573"
$ws.Range("C20").Value = [double]"0.0002533237224135188"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"0.07146556676033214"
$ws.Range("C21").Value = [double]"0.8784604219864331"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.555650872940484"
$ws.Range("C22").Value = [double]"0.006407309883128588"

# ---- Sheet 9: summ1 -> summ12 ----
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ12"
$ws.Rows("5:5").Delete()
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = [double]"-0.374315702615049"
$ws.Range("C2").Value = [double]"0.6343626497921968"
$ws.Range("A3").Value = "HHType_simp[T.MultiAdult_Kids]"
$ws.Range("B3").Value = [double]"0.2483227052521567"
$ws.Range("C3").Value = [double]"0.06554227857143982"
$ws.Range("A4").Value = "HHType_simp[T.Single_Female]"
$ws.Range("B4").Value = [double]"-0.5250599207462631"
$ws.Range("C4").Value = [double]"1.560788417449264e-06"
$ws.Range("A5").Value = "HHType_simp[T.Single_Male]"
$ws.Range("B5").Value = [double]"-0.3122485854078308"
$ws.Range("C5").Value = [double]"0.005483514300701871"
$ws.Range("A6").Value = "HHType_simp[T.Single_Parent]"
$ws.Range("B6").Value = [double]"-0.3248733419457782"
$ws.Range("C6").Value = [double]"0.0102367201536421"
$ws.Range("A7").Value = "HHSize"
$ws.Range("B7").Value = [double]"0.02752115045091767"
$ws.Range("C7").Value = [double]"0.5926883805231311"
$ws.Range("A8").Value = "IncomeDetailed_Numeric"
$ws.Range("B8").Value = [double]"0.0006257460278955607"
$ws.Range("C8").Value = [double]"1.180411082138264e-96"
$ws.Range("A9").Value = "maxAgeHH"
$ws.Range("B9").Value = [double]"0.01359773534731856"
$ws.Range("C9").Value = [double]"6.412953573675815e-07"
$ws.Range("A10").Value = "UniversityEducation"
$ws.Range("B10").Value = [double]"0.3327842857239968"
$ws.Range("C10").Value = [double]"2.823190363067429e-06"
$ws.Range("A11").Value = "InEmployment"
$ws.Range("B11").Value = [double]"0.2707570964657362"
$ws.Range("C11").Value = [double]"0.01476652368774869"
$ws.Range("A12").Value = "AllRetired"
$ws.Range("B12").Value = [double]"0.3757891915426358"
$ws.Range("C12").Value = [double]"0.008664729782130834"
$ws.Range("A13").Value = "UrbPopDensity"
$ws.Range("B13").Value = [double]"-2.739794801831464e-05"
$ws.Range("C13").Value = [double]"2.381207660919387e-05"
$ws.Range("A14").Value = "UrbBuildDensity"
$ws.Range("B14").Value = [double]"-1.40549439929387e-08"
$ws.Range("C14").Value = [double]"0.2041552722673118"
$ws.Range("A15").Value = "DistSubcenter"
$ws.Range("B15").Value = [double]"-0.006689738138297396"
$ws.Range("C15").Value = [double]"0.7136984945326275"
$ws.Range("A16").Value = "DistCenter"
$ws.Range("B16").Value = [double]"0.04870817744020787"
$ws.Range("C16").Value = [double]"0.001285946204447911"
$ws.Range("A17").Value = "bike_lane_share"
$ws.Range("B17").Value = [double]"-2.128770832126347"
$ws.Range("C17").Value = [double]"3.877359991486935e-09"
$ws.Range("A18").Value = "IntersecDensity"
$ws.Range("B18").Value = [double]"-0.006408621899093375"
$ws.Range("C18").Value = [double]"0.05827063499927885"
$ws.Range("A19").Value = "StreetLength"
$ws.Range("B19").Value = [double]"-0.002978781058516087"
$ws.Range("C19").Value = [double]"0.5246909703321185"
$ws.Range("A20").Value = "LU_UrbFab"
$ws.Range("B20").Value = [double]"1.332657994228911"
$ws.Range("C20").Value = [double]"0.0003787510569027891"
$ws.Range("A21").Value = "LU_Comm"
$ws.Range("B21").Value = [double]"0.3854121518600655"
$ws.Range("C21").Value = [double]"0.4140526252764645"
$ws.Range("A22").Value = "LU_Urban"
$ws.Range("B22").Value = [double]"-1.32674571556488"
$ws.Range("C22").Value = [double]"0.0187140850255425"

